$wb = $excel.ActiveWorkbook

$wsLogin = $wb.Worksheets.Item("Login")
$wsProject = $wb.Worksheets.Item("Add Project")

# --- Login sheet: new "error message" column (C) ---------------------------
$wsLogin.Range("C1").Value = "error message"
$wsLogin.Range("C2").Value = "Missing required parameter USERNAME"
$wsLogin.Range("C3").Value = "Missing required parameter USERNAME"
$wsLogin.Range("C4").Value = "Incorrect username or password."
$wsLogin.Range("C5").Value = "Incorrect username or password."
$wsLogin.Range("C6").Value = "Incorrect username or password."

# Wrap the long error-message cells so they read nicely in the narrower rows
$wsLogin.Range("C2:C3").WrapText = $true

# Widen column C to fit the new text instead of the old bestFit width
$wsLogin.Columns.Item(3).ColumnWidth = 43.3

# Page setup for the Login sheet (portrait, letter-ish A4/Letter paper id 9)
$wsLogin.PageSetup.PaperSize = 9
$wsLogin.PageSetup.Orientation = 1

# --- Active sheet / selection bookkeeping -----------------------------------
# Project sheet had the tab selection before; move it back to Login and park
# the cursor on C8 (first empty row under the new column) to match the saved
# view state captured in the workbook.
[void]$wsProject.Range("I4").Select()
[void]$wsLogin.Select()
[void]$wsLogin.Range("C8").Select()
